$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.399.20'
$ws.Range("E2").Value = '  -0.32%  '

# Row 3
$ws.Range("D3").Value = '1.823.15'
$ws.Range("E3").Value = '  -0.64%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '315.01'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -1.11%  '

# Row 6
$ws.Range("E6").Value = '  -0.02%  '

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.5134'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -3.48%  '

# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.3926'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -3.70%  '

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.07667'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +1.32%  '

# Row 10
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '1.109'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -0.20%  '

# Row 11
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '41.57'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -0.72%  '

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '21.01'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +0.76%  '

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '6.269'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -1.01%  '

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +0.03%  '

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '7.489'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -1.89%  '

# Row 16
$ws.Range("D16").Value = '1.824.26'
$ws.Range("E16").Value = '  -0.77%  '

# Row 17
$ws.Range("E17").Value = '  +3.99%  '

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '0.00001095'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +2.03%  '

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.06664'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +1.03%  '

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '17.69'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +0.60%  '

# Row 21
$ws.Range("E21").Value = '  +0.08%  '

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '6.121'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +0.77%  '

# Row 23
$ws.Range("D23").Value = '28.419.16'
$ws.Range("E23").Value = '  -0.29%  '

# Row 24
$ws.Range("E24").Value = '  -1.50%  '

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.254'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +6.64%  '

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '20.78'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +0.97%  '

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '156.41'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -0.32%  '

# Row 28
$ws.Range("D28").Value = '2.034.84'
$ws.Range("E28").Value = '  -0.64%  '

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '2.391'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -2.61%  '

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '124.11'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +0.13%  '

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '1.110'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -1.47%  '

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '0.1091'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -0.34%  '

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '5.649'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -0.80%  '

# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '3.657'
$cell.Style = "Normal"

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.07119'
$cell.Style = "Normal"

# Row 36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.2210'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -2.78%  '

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '0.02326'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -1.01%  '

# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '5.169'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -2.03%  '

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '8.789'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -0.75%  '

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.6256'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -0.46%  '

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '11.21'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -1.34%  '

# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '1.169'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -2.03%  '

# Row 43
$ws.Range("E43").Value = '  +0.01%  '

# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '1.392'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -1.46%  '

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '13.32'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -0.73%  '

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '3.716'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -0.06%  '

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.5882'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +0.35%  '

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '124.58'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -1.00%  '

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '1.981'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -0.52%  '

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '1.195'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +0.01%  '

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.06897'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -0.15%  '
